$wb = $excel.ActiveWorkbook

# --- TS_Defs sheet: change cell Q6 from "p,t" to "t" ---
$tsDefs = $wb.Worksheets.Item("TS_Defs")
$tsDefs.Range("Q6").Value = "t"

# --- process map sheet: add two new rows (24 and 25) ---
$procMap = $wb.Worksheets.Item("process map")
$procMap.Range("A24").Value = "old_new"
$procMap.Range("B24").Value = "*"
$procMap.Range("C24").Value = "new"

$procMap.Range("A25").Value = "old_new"
$procMap.Range("B25").Value = "ep*"
$procMap.Range("C25").Value = "old"

# --- selection / active sheet updates to match final state ---
$tsDefs.Range("A6").Select()
$procMap.Select()
$procMap.Range("C26").Select()
